# TC07_C3DC_phs000463_DiseasePhase-InitialDiag.xlsx
# "Updated remaining queries for C3DC"
#
# The SQL text stored in the TabQuery / StatQuery cells referenced the old
# join-key names (std.id / prt.id / "study.id" / "participant.id"). Update
# every occurrence to the new, fully-qualified key names
# (std.study_id / prt.participant_id / "study.study_id" / "participant.participant_id").

function Fix-Query([string]$sql) {
    $sql = $sql.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $sql = $sql.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $sql = $sql.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $sql = $sql.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $sql = $sql.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $sql = $sql.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $sql
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 = StatQuery, B2:B7 = TabQuery (dbGaP/Participants/Diagnosis/Treatment/
# TreatmentResp/Survival) -- every query in the sheet joins on the same
# study/participant keys, so all seven cells get the same fix.
$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $queryCells) {
    $cell = $ws.Range($addr)
    $current = $cell.Value()
    $cell.Value = Fix-Query $current
}

# Selection moved from C7 back to B2 (and the saved top-left scroll anchor
# is cleared along with it).
[void]$ws.Range("B2").Select()

# Column C is no longer auto "best fit" -- it now has an explicit width.
$ws.Columns.Item(3).ColumnWidth = 69.33
